# Ran code for averaged intensities on spiral schemes.
#
# The "Gaussian-Quadrature" averaging scheme is moved up in the table
# (right after "Ring Perpendicular to TD") and three new spiral-based
# averaging schemes are appended, each with a full row of averaged
# intensity values (1 for every HKL peak column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters for the 14 HKL-peak columns (C..P).
$peakCols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P")

# Final ordering of averaging-scheme row labels for rows 3-19.
$schemeNames = @(
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "Gaussian-Quadrature",
    "Spiral-90deg-10rot-5space",
    "Spiral-90deg-15rot-5space",
    "Spiral-90deg-10rot-3space",
    "NoRotation-tilt60deg",
    "Rotation-NoTilt",
    "Rotation-60detTilt",
    "HexGrid-90degTilt5degRes",
    "HexGrid-90degTilt22p5degRes",
    "HexGrid-60degTilt5degRes"
)

for ($i = 0; $i -lt $schemeNames.Length; $i++) {
    $row = $i + 3
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $schemeNames[$i]
    foreach ($col in $peakCols) {
        $ws.Range($col + $row).Value = 1
    }
}

# The three brand-new rows (17-19) need the same "index" cell formatting
# (bold, centered, bordered) that the existing index column already uses
# for rows 3-16 - copy it across instead of restyling from scratch so we
# don't introduce any extra/unused style definitions.
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
